# Update gh-pages output data: refresh "想去人数" (F) counts (and one "最低票价" G
# correction) on both the "展览" sheet and the aggregated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find each target row by matching the event name in column C, then update
    # the "想去人数" value in column F (and, for one row, the "最低票价" in G).
    $updates = @{
        "南昌·SuperComic动漫游戏博览会" = @{ F = 3809 }
        "宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华" = @{ F = 97 }
        "赣州·第四届赣州半夏动漫展" = @{ F = 936 }
        "抚州·临次元08·盛夏动漫狂欢节" = @{ F = 97 }
        "南昌·萌卡动漫展" = @{ F = 3352 }
        "江西·次元星河动漫游戏嘉年华" = @{ F = 5686 }
        "南昌·幻梦境国际动漫游戏嘉年华1th" = @{ F = 3327 }
        "吉安·COMIC LIFE周年庆典" = @{ F = 346 }
        "景德镇·第十五届瓷都ACG动漫游戏博览会" = @{ F = 2433 }
        "景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票" = @{ F = 568 }
        "萍乡·AU9夏至国漫展" = @{ F = 115 }
        "赣州·第一届环梦动漫游戏嘉年华" = @{ F = 191; G = 50 }
        "九江·第一届异次元动漫嘉年华" = @{ F = 343 }
        "上饶·囧喵喵国风动漫展" = @{ F = 111 }
        "赣州·第二届异次元动漫嘉年华" = @{ F = 886 }
        "南昌·W·MEETING动漫游戏盛典" = @{ F = 13 }
        "吉安·WF无线次元新星动漫博览会" = @{ F = 31 }
        "上饶·次元重现夏日嘉年华" = @{ F = 59 }
        "南昌·第四届龙年动漫展——暑假最后的狂欢" = @{ F = 541 }
    }

    $dims = $ws.UsedRange
    $lastRow = $dims.Rows.Count

    for ($r = 1; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($null -eq $name) { continue }
        if ($updates.ContainsKey($name)) {
            $upd = $updates[$name]
            if ($upd.ContainsKey("F")) {
                $ws.Cells.Item($r, 6).Value = $upd["F"]
            }
            if ($upd.ContainsKey("G")) {
                $ws.Cells.Item($r, 7).Value = $upd["G"]
            }
        }
    }
}
